$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the policy number in E2 (stored as text, quote-prefixed)
$ws.Range("E2").Value = "'12112001753"

# Move the active selection to E3 (matches the saved view state)
$ws.Range("E3").Select()
